# Natmi following Dr Hou advice: refresh the Ptn-Ptprb LR-pairs sheet
# with the updated cluster set (adds an "ECs" sending/target cluster)
# and recomputed specificity/expression statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-10 (columns A-T), row-major.
$data = @(
  ,@("ECs", "Ptn", "Ptprb", "ECs", 1, 0.3333333333333333, 0.3206776666666667, 0.962033, 0.03001977461414601, 0.03001977461414601, 3, 1, 141.7673063333333, 425.301919, 0.9653621431809045, 0.9653621431809047, 45.46160900459189, 409.154481041327, 0.0289799539593197, 0.02897995395931971)
  ,@("ECs", "Ptn", "Ptprb", "FAPs", 1, 0.3333333333333333, 0.3206776666666667, 0.962033, 0.03001977461414601, 0.03001977461414601, 3, 1, 3.731893333333333, 11.19568, 0.02541226633677049, 0.02541226633677049, 1.196734846382222, 10.77061361744, 0.0007628705078644999, 0.0007628705078645)
  ,@("ECs", "Ptn", "Ptprb", "sCs", 1, 0.3333333333333333, 0.3206776666666667, 0.962033, 0.03001977461414601, 0.03001977461414601, 3, 1, 1.354815, 4.064445, 0.009225590482324892, 0.009225590482324892, 0.4344589129650001, 3.910130216685, 0.0002769501469618039, 0.0002769501469618039)
  ,@("FAPs", "Ptn", "Ptprb", "ECs", 3, 1, 5.752274333333333, 17.256823, 0.5384908178993973, 0.5384908178993975, 3, 1, 141.7673063333333, 425.301919, 0.9653621431809045, 0.9653621431809047, 815.4844375270374, 7339.359937743337, 0.5198386500506004, 0.5198386500506006)
  ,@("FAPs", "Ptn", "Ptprb", "FAPs", 3, 1, 5.752274333333333, 17.256823, 0.5384908178993973, 0.5384908178993975, 3, 1, 3.731893333333333, 11.19568, 0.02541226633677049, 0.02541226633677049, 21.46687423607111, 193.20186812464, 0.01368427208436486, 0.01368427208436487)
  ,@("FAPs", "Ptn", "Ptprb", "sCs", 3, 1, 5.752274333333333, 17.256823, 0.5384908178993973, 0.5384908178993975, 3, 1, 1.354815, 4.064445, 0.009225590482324892, 0.009225590482324892, 7.793267550915001, 70.139407958235, 0.004967895764432027, 0.004967895764432028)
  ,@("sCs", "Ptn", "Ptprb", "ECs", 3, 1, 4.609262333333334, 13.827787, 0.4314894074864565, 0.4314894074864565, 3, 1, 141.7673063333333, 425.301919, 0.9653621431809045, 0.9653621431809047, 653.4427051803615, 5880.984346623253, 0.4165435391709843, 0.4165435391709844)
  ,@("sCs", "Ptn", "Ptprb", "FAPs", 3, 1, 4.609262333333334, 13.827787, 0.4314894074864565, 0.4314894074864565, 3, 1, 3.731893333333333, 11.19568, 0.02541226633677049, 0.02541226633677049, 17.20127537335111, 154.81147836016, 0.01096512374454112, 0.01096512374454113)
  ,@("sCs", "Ptn", "Ptprb", "sCs", 3, 1, 4.609262333333334, 13.827787, 0.4314894074864565, 0.4314894074864565, 3, 1, 1.354815, 4.064445, 0.009225590482324892, 0.009225590482324892, 6.244697748135001, 56.202279733215, 0.003980744570931061, 0.003980744570931061)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $data[$i]
  $r = $startRow + $i
  for ($c = 1; $c -le $row.Count; $c++) {
    $ws.Cells.Item($r, $c).Value = $row[$c - 1]
  }
}
